# Bugfix: complete data by using values both from codelist and data (#80)
# Close #79
#
# The "Codelists" sheet is missing a code-list row for age "99" (a
# catch-all/missing-value code). Insert a new row 103 with that entry,
# which pushes the existing rows 103-111 down to 104-112 (values/styles
# unchanged), then fill the new row with the age/99 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Codelists")

# Insert a new blank row at 103 (shifts old rows 103:111 -> 104:112).
$ws.Rows("103:103").Insert() | Out-Null

# Populate the new row: age codelist entry "99" (position/en/da/kl/sortorder).
$ws.Range("A103").Value = "age"
$ws.Range("B103").Value = 99
$ws.Range("C103").Value = 99
$ws.Range("D103").Value = 99
$ws.Range("E103").Value = 99
$ws.Range("F103").Value = 99

# Make "Codelists" the active sheet and put the selection on the new row's
# precision column, matching the saved UI state in the workbook.
$ws.Activate()
$ws.Range("G103").Select()
